$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text/string assignments (non-numeric-looking values stay text automatically)
$ws.Range("D2").Value = "64.088.10"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.068.57"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("E6").Value = "  +4.64%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.064.85"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "3.571.39"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "63.995.41"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "3.055.58"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  +4.81%  "
$ws.Range("E24").Value = "  +7.70%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("E29").Value = "  +3.60%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("E38").Value = "  +15.66%  "
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("D41").Value = "2.958.10"
$ws.Range("E41").Value = "  -6.09%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E46").Value = "  +4.20%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("E51").Value = "  +0.26%  "

# Numeric-looking Price values: force text via NumberFormat, then restore default formatting
$ws.Range("D4:D7").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D5").Value = "557.58"
$ws.Range("D6").Value = "145.68"
$ws.Range("D7").Value = "0.999"
$ws.Range("Z1").Copy()
$ws.Range("D4:D7").PasteSpecial(-4122)

$ws.Range("D9:D14").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D10").Value = "0.152"
$ws.Range("D11").Value = "6.30"
$ws.Range("D12").Value = "0.470"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D14").Value = "35.22"
$ws.Range("Z1").Copy()
$ws.Range("D9:D14").PasteSpecial(-4122)

$ws.Range("D19:D31").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("D20").Value = "476.36"
$ws.Range("D21").Value = "13.93"
$ws.Range("D22").Value = "0.675"
$ws.Range("D23").Value = "7.54"
$ws.Range("D24").Value = "13.52"
$ws.Range("D25").Value = "81.63"
$ws.Range("D26").Value = "1.00"
$ws.Range("D27").Value = "2.79"
$ws.Range("D28").Value = "8.09"
$ws.Range("D29").Value = "2.06"
$ws.Range("D30").Value = "0.997"
$ws.Range("D31").Value = "26.17"
$ws.Range("Z1").Copy()
$ws.Range("D19:D31").PasteSpecial(-4122)

$ws.Range("D33:D40").NumberFormat = "@"
$ws.Range("D33").Value = "2.49"
$ws.Range("D34").Value = "5.59"
$ws.Range("D35").Value = "6.16"
$ws.Range("D36").Value = "54.85"
$ws.Range("D37").Value = "462.67"
$ws.Range("D38").Value = "3.00"
$ws.Range("D39").Value = "0.0831"
$ws.Range("D40").Value = "0.0405"
$ws.Range("Z1").Copy()
$ws.Range("D33:D40").PasteSpecial(-4122)

$ws.Range("D42:D49").NumberFormat = "@"
$ws.Range("D42").Value = "8.28"
$ws.Range("D43").Value = "0.114"
$ws.Range("D44").Value = "28.08"
$ws.Range("D45").Value = "0.260"
$ws.Range("D46").Value = "2.14"
$ws.Range("D47").Value = "1.00"
$ws.Range("D48").Value = "0.112"
$ws.Range("D49").Value = "119.72"
$ws.Range("Z1").Copy()
$ws.Range("D42:D49").PasteSpecial(-4122)

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("Z1").Copy()
$ws.Range("D51").PasteSpecial(-4122)

$excel.CutCopyMode = 0
